$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Distinguish CHE1 data files from CHE2 by renaming the relative path entries
# for the CHE1 study rows (ageband + region) to include the "1" suffix.
$ws.Range("C8").Value = "data/derived/CHE/CHE1_agebands.RDS"
$ws.Range("C9").Value = "data/derived/CHE/CHE1_region.RDS"

# Update the active cell selection on Sheet1 to match the saved workbook state
$ws.Range("C10").Select()
